# The source export was re-run a day later (2024-08-15 09:46:16 -> 2024-08-16
# 10:14:25), so the sheet name carrying that timestamp changes, every row's
# "Dt. Referencia" (column G) advances by one day (Excel serial 45519 -> 45520),
# and two accounts (rows 231 & 232) picked up updated balances for the new day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to match the new export timestamp.
$ws.Name = "IClientBalance-20240816-101425-"

# Bump the reference date (column G) for every data row (2..274) by one day.
$dataRowCount = 274
for ($r = 2; $r -le $dataRowCount; $r++) {
    $ws.Cells.Item($r, 7).Value = 45520
}

# Row 231 ("Saldo Previsto" / "Vl. Total") now reports 39858.22 instead of 28881.77.
$ws.Cells.Item(231, 5).Value = 39858.22
$ws.Cells.Item(231, 8).Value = 39858.22

# Row 232 ("Saldo Previsto" / "Vl. Total") now reports 34315.74 instead of 1370.9.
$ws.Cells.Item(232, 5).Value = 34315.74
$ws.Cells.Item(232, 8).Value = 34315.74

# Reset the lingering cell selection left over from the previous save (was B2).
$ws.Range("A1").Select() | Out-Null
